# Update cryptocurrency price/volume figures with the latest snapshot values.
# Values are stored as plain text (not numbers) in the source sheet, so each
# target cell is temporarily switched to the "@" text number format before the
# new value is written (otherwise Excel auto-converts numeric-looking text into
# a real number/percentage). The style is then reset back to "Normal" so the
# cell's formatting matches the rest of the untouched sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "D2" = "325.64"
    "E2" = "-1.46%"
    "E3" = "-1.13%"
    "D4" = "5.725"
    "E4" = "6.49%"
    "D5" = "0.08036"
    "E5" = "-0.89%"
    "D6" = "2.031"
    "E6" = "5.29%"
    "D7" = "4.496"
    "E7" = "-0.70%"
    "D8" = "8.638"
    "E8" = "-0.17%"
    "D9" = "2.930"
    "E9" = "-1.58%"
    "D10" = "0.9228"
    "E10" = "-2.06%"
    "D11" = "0.1256"
    "E11" = "-7.64%"
    "D12" = "0.1955"
    "E12" = "-0.85%"
    "D13" = "8.754"
    "E13" = "20.73%"
    "D14" = "0.09281"
    "E14" = "0.09%"
    "D15" = "0.03570"
    "E15" = "0.14%"
    "E16" = "9.40%"
    "E17" = "-2.58%"
    "D18" = "0.006161"
    "E18" = "-0.52%"
    "D19" = "3.363"
    "E19" = "-0.10%"
    "D20" = "0.3480"
    "E20" = "-1.16%"
    "E21" = "1.47%"
    "D22" = "0.2666"
    "E22" = "4.16%"
    "D23" = "0.04406"
    "E23" = "-0.65%"
    "D24" = "0.001261"
    "E24" = "3.24%"
    "D25" = "0.004608"
    "E25" = "7.78%"
    "D26" = "0.0001191"
    "E26" = "-0.75%"
    "D39" = "0.02494"
    "E39" = "0.14%"
    "D40" = "0.05315"
    "E40" = "1.79%"
    "D41" = "0.007478"
    "E41" = "-0.91%"
    "D42" = "0.009911"
    "E42" = "8.50%"
    "E43" = "-1.59%"
    "D45" = "0.01151"
    "E45" = "6.71%"
    "D46" = "0.00006682"
    "E46" = "1.54%"
    "D47" = "0.00000000750"
    "D48" = "0.003040"
    "E48" = "-9.09%"
    "E49" = "-5.03%"
    "D50" = "0.00002101"
    "D51" = "0.0002001"
}

foreach ($ref in $newValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$ref]
    $cell.Style = "Normal"
}
